# Fill in column I (Reconcile Action) with the resolved values on the
# "PUR_unresolved_case" sheet (replacing the blanket "unresolved_case"
# placeholder with the actual reconciled land-use classification).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PUR_unresolved_case")

$values = @{
    2 = "Hutan Lindung"
    3 = "Hutan Lindung"
    4 = "Hutan_Adat_Bungo_F"
    5 = "HTI_Bungo_F"
    6 = "Tambang_Bungo_F"
    7 = "Hutan Produksi"
    8 = "Hutan Produksi"
    9 = "HTI_Bungo_F"
    10 = "Hutan Produksi"
    11 = "HTI_Bungo_F"
    12 = "Perkebunan"
    13 = "HTI_Bungo_F"
    14 = "Tambang_Bungo_F"
    15 = "Perkebunan"
    16 = "HTI_Bungo_F"
    17 = "Perkebunan"
    18 = "Pertanian Lahan Basah"
    19 = "Pertanian Lahan Basah"
    20 = "Pertanian Lahan Basah"
    21 = "Pertanian Lahan Basah"
    22 = "HTI_Bungo_F"
    23 = "Tambang_Bungo_F"
    24 = "Pertanian Lahan Kering"
    25 = "Pertanian Lahan Kering"
    26 = "Pertanian Lahan Kering"
    27 = "Sempadan Sungai"
    28 = "Sempadan Sungai"
    29 = "Sempadan Sungai"
    30 = "Sempadan Sungai"
}

foreach ($row in $values.Keys) {
    $ws.Range("I$row").Value = $values[$row]
}

$ws.Columns.Item(7).ColumnWidth = 23.833333333333332
$ws.Range("K27").Select()
